$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted as the most-recent row for
# "Cilantro" at "Terminal La Palmera de La Serena" (row 205), pushing the
# previously-existing rows 205-215 down to 206-216.
$ws.Rows.Item(205).Insert()

$ws.Cells.Item(205, 1).Value = 8
$ws.Cells.Item(205, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(205, 3).Value = "Coquimbo"
$ws.Cells.Item(205, 4).Value = 45008
$ws.Cells.Item(205, 5).Value = 4
$ws.Cells.Item(205, 6).Value = 100112040
$ws.Cells.Item(205, 7).Value = "Cilantro"
$ws.Cells.Item(205, 8).Value = "Sin especificar"
$ws.Cells.Item(205, 9).Value = "Primera"
$ws.Cells.Item(205, 10).Value = 2400
$ws.Cells.Item(205, 11).Value = 1800
$ws.Cells.Item(205, 12).Value = 2000
$ws.Cells.Item(205, 13).Value = 1900
$ws.Cells.Item(205, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(205, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(205, 16).Value = 1267
$ws.Cells.Item(205, 17).Value = 1.5
$ws.Cells.Item(205, 18).Value = "Hortaliza"
